$d = $word.ActiveDocument

# The document currently ends with a ListParagraph ("Added a placeholder sprite")
# that carries the trailing _GoBack bookmark. We append two new list items
# after it, matching its paragraph/run formatting (style, numbering, spacing,
# fonts), and give them the new text.

$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$d.Paragraphs.Last.Range.Text = "Added a placeholder bg"

$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$d.Paragraphs.Last.Range.Text = "Added pseudocode for scroll() method"
